$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Find-RowByAccount($acc) {
    $lastRow = $ws.UsedRange.Rows.Count
    for ($r = 1; $r -le $lastRow; $r++) {
        $a = $ws.Cells.Item($r, 1).Text
        if ($a -eq $acc) {
            return $r
        }
    }
    return -1
}

function Remove-AccountRow($acc) {
    $r = Find-RowByAccount $acc
    if ($r -gt 0) {
        $ws.Cells.Item($r, 1).EntireRow.Delete()
    }
}

function Insert-AccountRowBefore($beforeAcc, $acc, $name, $val) {
    $r = Find-RowByAccount $beforeAcc
    $ws.Cells.Item($r, 1).EntireRow.Insert()
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $acc
    $ws.Cells.Item($r, 2).Value = $name
    $ws.Cells.Item($r, 3).Value = $val
}

# Remove rows that no longer exist in the export
Remove-AccountRow "004468717"   # HELOISA 23079
Remove-AccountRow "005547703"   # SILVIA 22296.58
Remove-AccountRow "004361159"   # HFR 1013.62
Remove-AccountRow "004451652"   # MATEUS 97.9

# MARINA's old-position row (account 003249855, balance 86.22) is removed; a
# new row for the same account is inserted further up with an updated balance.
Remove-AccountRow "003249855"

# Insert the new CARLOS row (827.8) right before GISELA (004322719)
Insert-AccountRowBefore "004322719" "005685353" "CARLOS" 827.8

# Insert MARINA's updated row (237.12) right before ANA (004467884)
Insert-AccountRowBefore "004467884" "003249855" "MARINA" 237.12
